$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that must stay plain text (e.g. numeric-looking price
# strings like "1.00" or "0.110") without Excel converting it into a real number.
# We briefly mark the cell as Text, set the value, then restore the "Normal" style
# so the saved file does not carry a stray number-format style index.
function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = '63.383.33'
$ws.Range("E2").Value = '  -7.35%  '

# Row 3
$ws.Range("D3").Value = '3.539.20'
$ws.Range("E3").Value = '  -2.85%  '

# Row 4
$ws.Range("E4").Value = '  +0.13%  '

# Row 5
Set-TextValue $ws.Range("D5") '390.31'
$ws.Range("E5").Value = '  -7.29%  '

# Row 6
Set-TextValue $ws.Range("D6") '123.31'
$ws.Range("E6").Value = '  -6.49%  '

# Row 7
$ws.Range("D7").Value = '3.534.58'
$ws.Range("E7").Value = '  -2.63%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.589'
$ws.Range("E8").Value = '  -11.27%  '

# Row 9
Set-TextValue $ws.Range("D9") '1.00'
$ws.Range("E9").Value = '  +0.09%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.681'
$ws.Range("E10").Value = '  -12.27%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.151'
$ws.Range("E11").Value = '  -23.75%  '

# Row 12
Set-TextValue $ws.Range("D12") '0.0000320'
$ws.Range("E12").Value = '  -25.06%  '

# Row 13
Set-TextValue $ws.Range("D13") '38.70'
$ws.Range("E13").Value = '  -8.90%  '

# Row 14
$ws.Range("D14").Value = '4.112.27'
$ws.Range("E14").Value = '  -2.15%  '

# Row 15
Set-TextValue $ws.Range("D15") '9.16'
$ws.Range("E15").Value = '  -7.66%  '

# Row 16
Set-TextValue $ws.Range("D16") '0.136'
$ws.Range("E16").Value = '  -2.85%  '

# Row 17
$ws.Range("D17").Value = '3.537.26'
$ws.Range("E17").Value = '  -2.00%  '

# Row 18
Set-TextValue $ws.Range("D18") '12.99'
$ws.Range("E18").Value = '  +3.91%  '

# Row 19
Set-TextValue $ws.Range("D19") '18.74'
$ws.Range("E19").Value = '  -7.43%  '

# Row 20
$ws.Range("D20").Value = '63.518.38'
$ws.Range("E20").Value = '  -7.07%  '

# Row 21
Set-TextValue $ws.Range("D21") '1.02'
$ws.Range("E21").Value = '  -9.97%  '

# Row 22
Set-TextValue $ws.Range("D22") '393.09'
$ws.Range("E22").Value = '  -14.73%  '

# Row 23
Set-TextValue $ws.Range("D23") '13.83'
$ws.Range("E23").Value = '  +3.16%  '

# Row 24
Set-TextValue $ws.Range("D24") '80.66'
$ws.Range("E24").Value = '  -9.96%  '

# Row 25
Set-TextValue $ws.Range("D25") '2.88'
$ws.Range("E25").Value = '  -6.94%  '

# Row 26
$ws.Range("E26").Value = '  +11.82%  '

# Row 27
Set-TextValue $ws.Range("D27") '33.67'
$ws.Range("E27").Value = '  -6.36%  '

# Row 28
Set-TextValue $ws.Range("D28") '2.98'
$ws.Range("E28").Value = '  -10.87%  '

# Row 29
Set-TextValue $ws.Range("D29") '8.70'
$ws.Range("E29").Value = '  -14.54%  '

# Row 30
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D30") '2.65'
$ws.Range("E30").Value = '  -4.77%  '

# Row 31
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range("D31") '11.78'
$ws.Range("E31").Value = '  -5.43%  '

# Row 32
Set-TextValue $ws.Range("D32") '0.110'
$ws.Range("E32").Value = '  -7.21%  '

# Row 33
Set-TextValue $ws.Range("D33") '6.79'
$ws.Range("E33").Value = '  -6.15%  '

# Row 34
$ws.Range("E34").Value = '  -5.06%  '

# Row 35
Set-TextValue $ws.Range("D35") '0.998'
$ws.Range("E35").Value = '  -0.06%  '

# Row 36
Set-TextValue $ws.Range("D36") '36.79'
$ws.Range("E36").Value = '  -8.38%  '

# Row 37
Set-TextValue $ws.Range("D37") '53.90'
$ws.Range("E37").Value = '  -3.99%  '

# Row 38
Set-TextValue $ws.Range("D38") '0.0439'
$ws.Range("E38").Value = '  -10.10%  '

# Row 39
Set-TextValue $ws.Range("D39") '0.999'
$ws.Range("E39").Value = '  -0.02%  '

# Row 40
Set-TextValue $ws.Range("D40") '2.71'
$ws.Range("E40").Value = '  +7.16%  '

# Row 41
Set-TextValue $ws.Range("D41") '0.131'
$ws.Range("E41").Value = '  -12.50%  '

# Row 42
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D42") '142.39'
$ws.Range("E42").Value = '  -4.31%  '

# Row 43
$ws.Range("B43").Value = 'ApeXProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws.Range("D43") '3.05'
$ws.Range("E43").Value = '  +14.55%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D44") '25.55'
$ws.Range("E44").Value = '  +17.86%  '

# Row 45
$ws.Range("B45").Value = 'PEPE'
$ws.Range("C45").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D45").Value = '0.0₃0599'
$ws.Range("E45").Value = '  -29.24%  '

# Row 46
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D46") '2.51'
$ws.Range("E46").Value = '  -6.05%  '

# Row 47
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D47") '1.96'
$ws.Range("E47").Value = '  -0.72%  '

# Row 48
$ws.Range("B48").Value = 'LidoDAOToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range("D48") '3.08'
$ws.Range("E48").Value = '  -6.18%  '

# Row 49
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D49") '4.07'
$ws.Range("E49").Value = '  -5.21%  '

# Row 50
Set-TextValue $ws.Range("D50") '2.66'
$ws.Range("E50").Value = '  -10.20%  '

# Row 51
Set-TextValue $ws.Range("D51") '0.274'
$ws.Range("E51").Value = '  -9.41%  '
